# Auto-generated edit script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 30000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 90000
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -90460

$ws.Range("H40").Value = 1945.4546
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1880
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1880
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2230

$ws.Range("H129").Value = 3405.349
$ws.Range("J129").Value = 3652.0688
$ws.Range("L129").Value = 10956.2064
$ws.Range("N129").Value = -20956.2064

$ws.Range("H135").Value = 1888.3572
$ws.Range("I135").Value = 1312.7576
$ws.Range("J135").Value = 3998.889
$ws.Range("K135").Value = 11814.8184
$ws.Range("L135").Value = 35990.001
$ws.Range("M135").Value = -9279.8184
$ws.Range("N135").Value = -41060.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7043363.5
$ws.Range("I61").Value = 8475524
$ws.Range("K61").Value = 8475524
$ws.Range("M61").Value = -8475312

$ws.Range("H122").Value = 10421.583
$ws.Range("I122").Value = 11705.9
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 35117.7
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -32667.7
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 6251862
$ws.Range("I132").Value = 7354537
$ws.Range("J132").Value = 3369
$ws.Range("K132").Value = 22063611
$ws.Range("L132").Value = 10107
$ws.Range("M132").Value = -22061081
$ws.Range("N132").Value = -15167

$ws.Range("H136").Value = 7043363.5
$ws.Range("I136").Value = 8475524
$ws.Range("K136").Value = 25426572
$ws.Range("M136").Value = -25424022

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2154.5715
$ws.Range("I134").Value = 1359.6586
$ws.Range("J134").Value = 6228.5
$ws.Range("K134").Value = 4078.9758
$ws.Range("L134").Value = 18685.5
$ws.Range("M134").Value = -1543.9758
$ws.Range("N134").Value = -23755.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37040700
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 37040700
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 37040700
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -37041290

$ws.Range("H34").Value = 37040700
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 37040700
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 37040700
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -37041104

$ws.Range("H58").Value = 1337.2941
$ws.Range("I58").Value = 492.42856
$ws.Range("J58").Value = 5280
$ws.Range("K58").Value = 492.42856
$ws.Range("L58").Value = 5280
$ws.Range("M58").Value = -289.42856
$ws.Range("N58").Value = -5686

$ws.Range("H60").Value = 8189.25
$ws.Range("I60").Value = 7500
$ws.Range("J60").Value = 8419
$ws.Range("K60").Value = 7500
$ws.Range("L60").Value = 8419
$ws.Range("M60").Value = -6989
$ws.Range("N60").Value = -9441

$ws.Range("H122").Value = 2320.4285
$ws.Range("I122").Value = 2293.8462
$ws.Range("K122").Value = 6881.5386
$ws.Range("M122").Value = -4431.5386

$ws.Range("H134").Value = 1633.174
$ws.Range("I134").Value = 1729.1177
$ws.Range("K134").Value = 5187.3531
$ws.Range("M134").Value = -2652.3531

$ws.Range("H136").Value = 1337.2941
$ws.Range("I136").Value = 492.42856
$ws.Range("J136").Value = 5280
$ws.Range("K136").Value = 1477.28568
$ws.Range("L136").Value = 15840
$ws.Range("M136").Value = 1072.71432
$ws.Range("N136").Value = -20940

$ws.Range("H137").Value = 48689
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 48689
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 48689
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -58889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 603928.8
$ws.Range("I2").Value = 92.30768999999999
$ws.Range("J2").Value = 1388916.2
$ws.Range("K2").Value = 553.84614
$ws.Range("L2").Value = 8333497.199999999
$ws.Range("M2").Value = -440.84614
$ws.Range("N2").Value = -8333723.199999999

$ws.Range("H6").Value = 108.71429
$ws.Range("J6").Value = 55
$ws.Range("L6").Value = 165
$ws.Range("N6").Value = -391

$ws.Range("H12").Value = 86.545456
$ws.Range("J12").Value = 83.55556
$ws.Range("L12").Value = 250.66668
$ws.Range("N12").Value = -596.66668

$ws.Range("H39").Value = 536.57574
$ws.Range("J39").Value = 536.57574
$ws.Range("L39").Value = 1609.72722
$ws.Range("N39").Value = -2197.72722

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").ClearContents()

$ws.Range("H122").Value = 676.9259
$ws.Range("I122").Value = 626.0769
$ws.Range("K122").Value = 5634.6921
$ws.Range("M122").Value = -3184.6921

$ws.Range("H125").Value = 3148
$ws.Range("J125").Value = 3640
$ws.Range("L125").Value = 10920
$ws.Range("N125").Value = -20760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2667.4854
$ws.Range("I132").Value = 2009.5714
$ws.Range("J132").Value = 4364.2104
$ws.Range("K132").Value = 6028.7142
$ws.Range("L132").Value = 13092.6312
$ws.Range("M132").Value = -3498.7142
$ws.Range("N132").Value = -18152.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1264.3529
$ws.Range("I22").Value = 650.1667
$ws.Range("J22").Value = 1599.3636
$ws.Range("K22").Value = 650.1667
$ws.Range("L22").Value = 1599.3636
$ws.Range("M22").Value = -355.1667
$ws.Range("N22").Value = -2189.3636

$ws.Range("H27").Value = 1264.3529
$ws.Range("I27").Value = 650.1667
$ws.Range("J27").Value = 1599.3636
$ws.Range("K27").Value = 650.1667
$ws.Range("L27").Value = 1599.3636
$ws.Range("M27").Value = -543.1667
$ws.Range("N27").Value = -1813.3636
